$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: D3 (6-channel Nexperia TVS array) replaced with D3,D4,D5
#     (2-channel Toshiba TVS array), quantity bumped from 1 to 3 ---
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = "D3,D4,D5"
$ws.Range("E12").Value = "TVS diode array (2 channel)"
$ws.Range("G12").Value = "DF3A5.6LFULFCT-ND"
$ws.Range("F12").Value = "Toshiba Semiconductor and Storage"
$ws.Range("H12").Value = "DF3A5.6LFU,LF"
$ws.Range("I12").Value = 0.29

# --- Row 34: Control board PCB now ordered from OshPark, cost split 3 ways ---
$ws.Range("F34").Value = "OshPark"
$ws.Range("I34").Formula = "=33.76126/3"

# --- Row 35: LED strip price updated, formula rewritten explicitly ---
$ws.Range("I35").Formula = "=24.0864/3"
$ws.Range("J35").Formula = "=I35*C35"

# --- Row 44: Power board PCB now ordered from OshPark, cost split 3 ways ---
$ws.Range("F44").Value = "OshPark"
$ws.Range("I44").Formula = "=33.763596/3"
$ws.Range("J44").Formula = "=I44*C44"

# --- Sheet view / selection state ---
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("F48").Select()

$wb.Windows.Item(1).ScrollColumn = 5
